# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F3 and the Status column (C) on each language
#    sheet).
# 2) Narrow the "Status" column(s) that used to hold that text:
#    Overview columns E & F, and column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "In Translation"
$ov.Range("F2").Value = "In Translation"
$ov.Range("E3").Value = "In Translation"
$ov.Range("F3").Value = "In Translation"

$ov.Columns.Item(5).ColumnWidth = 12.5
$ov.Columns.Item(6).ColumnWidth = 12.5

# ---- zh-cn sheet ------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "In Translation"
$zh.Range("C3").Value = "In Translation"

$zh.Columns.Item(3).ColumnWidth = 12.5

# ---- de-de sheet ------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "In Translation"
$de.Range("C3").Value = "In Translation"

$de.Columns.Item(3).ColumnWidth = 12.5
